# Generate Report for Handback
#
# The localization file "0d2067ae-2ca4-448f-bc8e-89192d7c768c" has now been
# handed back (previously it was only "Ready for handoff"). Regenerate the
# report: flip row2/row3 ordering on the per-language detail sheets (and
# update the Overview sheet's status/date values + hyperlink display text)
# so that 0d2067ae's "Handed back" row leads, with a fresh handback
# timestamp, while e3b61664 keeps its already-known handback info.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Cell values: column A keeps its original file-name text per row, but the
# status/date columns (B/C/D) for row 3 now reflect a handback (same text
# used in row 2, since both are now handed back in sync with en-US).
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-23 18:57:33"

# Hyperlinks: the link targets (A2 -> e3b61664, A3 -> 0d2067ae) stay put,
# but the displayed text swaps.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/383e64ab3832172e41ec1ea1a8a090e3e07d53a3/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-23 18:57:29"
$ws2.Range("F2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws2.Range("G2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-23 18:58:05"
$ws2.Range("J2").Value = "Include"

$ws2.Range("A3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-23 18:55:47"
$ws2.Range("F3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws2.Range("G3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-23 18:56:33"
$ws2.Range("J3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f75b98accff473613980d6d00c49f8e386eb314/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/30fdea63b4e1cf4e0d65702035e6c8e5d3bb2d6f/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3436c8a9e40370fec3686ff289c357d4498874c7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/383e64ab3832172e41ec1ea1a8a090e3e07d53a3/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9186b2d2de244ed49d5a5dd163c7bb9bea89d516/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/30fdea63b4e1cf4e0d65702035e6c8e5d3bb2d6f/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3436c8a9e40370fec3686ff289c357d4498874c7/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.zh-cn.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-23 18:57:33"
$ws3.Range("F2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md"
$ws3.Range("G2").Value = "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-23 18:58:14"
$ws3.Range("J2").Value = "Include"

$ws3.Range("A3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-23 18:55:52"
$ws3.Range("F3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md"
$ws3.Range("G3").Value = "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-23 18:56:45"
$ws3.Range("J3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/91dfa5fecd8eb92d70a2f5fb7a12f3cf15aa0586/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cfd74deb0f2683a3fdf3da7265b4c3859c833f76/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4de9de7fdcccb267520dd644a0c0017c94b76d54/e2e/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7dfb35f40aea915747a2a1d008ebcfacd58ad633/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf", "", "", "0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/383e64ab3832172e41ec1ea1a8a090e3e07d53a3/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb257bdc2151c62b5b942c6ecc95f9dac6083d71/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4de9de7fdcccb267520dd644a0c0017c94b76d54/e2e/0d2067ae-2ca4-448f-bc8e-89192d7c768c.md", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7dfb35f40aea915747a2a1d008ebcfacd58ad633/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0d2067ae-2ca4-448f-bc8e-89192d7c768c.fdb0e26b708f7757c927665c32014118d81c9a82.de-de.xlf", "", "", "e3b61664-96dc-4ab8-bc89-9c0d7fefc835.5b909c6e9d6b4335b3af5e7cc35d338fdda5ef7d.de-de.xlf")
